$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '26.003.33'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E2").Value = '  +0.44%  '
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '1.642.08'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E3").Value = '  +0.11%  '
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.003'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E4").Value = '  +0.05%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '214.88'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E5").Value = '  -0.14%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.5093'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E6").Value = '  +1.23%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '1.003'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E7").Value = '  -0.30%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.2566'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E8").Value = '  -0.25%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.06364'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E9").Value = '  -0.34%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '19.57'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E10").Value = '  +0.20%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.07761'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E11").Value = '  -0.45%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '1.665.37'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E12").Value = '  +0.97%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '4.283'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E13").Value = '  +0.15%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '0.5443'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E14").Value = '  +0.38%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.0₅7735'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E15").Value = '  -1.85%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '64.21'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E16").Value = '  -0.75%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '26.037.22'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E17").Value = '  +0.19%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '1.003'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E18").Value = '  -0.43%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '196.89'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E19").Value = '  -0.41%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '4.420'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E20").Value = '  +0.79%  '
$ws.Range("E21").Value = '  -0.09%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '6.031'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E22").Value = '  +0.92%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '1.005'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E23").Value = '  -0.27%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '1.868'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E24").Value = '  -0.62%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '141.65'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E25").Value = '  +1.17%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '0.1189'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E26").Value = '  +4.36%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '6.824'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E27").Value = '  -0.30%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '15.57'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E28").Value = '  -0.80%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '1.236'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E29").Value = '  -0.54%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '0.04860'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E30").Value = '  -0.01%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '3.255'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E31").Value = '  -0.24%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '3.166'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E32").Value = '  -0.84%  '
$ws.Range("E33").Value = '  -0.55%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '2.367'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E34").Value = '  -0.21%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.8968'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E35").Value = '  +0.86%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '2.580'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E36").Value = '  -1.22%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '1.141.01'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E37").Value = '  +0.74%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.5454'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E38").Value = '  -1.47%  '
$ws.Range("E39").Value = '  -0.06%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '1.003'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E40").Value = '  -0.44%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '2.524'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E41").Value = '  -1.87%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.0₈128'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E42").Value = '  +8.30%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.8103'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E43").Value = '  -0.87%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '99.27'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E44").Value = '  -0.29%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '5.401'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E45").Value = '  -5.20%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '1.780.05'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E46").Value = '  +0.22%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.4533'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E47").Value = '  -0.03%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '54.89'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E48").Value = '  -0.60%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.9997'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E49").Value = '  -1.06%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.05058'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E50").Value = '  -0.49%  '
$ws.Range("E51").Value = '  -0.59%  '
